# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The account-statement table previously held a single worker
# (MARTHA LIGIA FONNEGRA GEORGE, CC 1143332154) with 5 overdue periods.
# It now holds two workers:
#   - EMIRO JOSE ORTIZ DURANGO (CC 15072746) with 7 overdue periods
#     (1812,1811,1810,1809,1808,1807,1806) at 36400 / 910000 each.
#   - MARTHA LIGIA FONNEGRA GEORGE (CC 1143332154) keeping her 5 periods
#     (1801,1712,1711,1710,1709) at 29600 / 840000 each.
# The summary header (total mora, worker count, period count) and the
# footer signature block move down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the worker table -------------------------------------------------
# Old layout: data rows 16-20 (one worker, 5 periods), closing-border style on
# row 20 (last row of the table).
# New layout: data rows 16-27 (two workers, 12 periods total), closing-border
# style on row 27 (new last row).
# Insert 7 blank rows at row 20; this pushes the old row 20 (and everything
# below it, including the footer) down by 7 rows, so it lands on row 27 and
# keeps its "closing" border style intact - which is exactly what the new
# last row needs.
$ws.Rows.Item(20).Resize(7).Insert()

# The freshly inserted rows 20-26 need the same "middle of the table" borders
# that rows 16-19 already use. Copy that formatting down from row 19.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J26").PasteSpecial(-4122)

# --- Populate the table content ---------------------------------------------
# Worker 1: EMIRO JOSE ORTIZ DURANGO - rows 16-22
$periodsEmiro = @("1812","1811","1810","1809","1808","1807","1806")
for ($i = 0; $i -lt $periodsEmiro.Length; $i++) {
  $r = 16 + $i
  $ws.Cells.Item($r, 2).Value = "CC"
  $ws.Cells.Item($r, 3).Value = "15072746"
  $ws.Cells.Item($r, 4).Value = "EMIRO JOSE ORTIZ DURANGO"
  $ws.Cells.Item($r, 5).Value = $periodsEmiro[$i]
  $ws.Cells.Item($r, 6).Value = 36400
  $ws.Cells.Item($r, 7).Value = 910000
}

# Worker 2: MARTHA LIGIA FONNEGRA GEORGE - rows 23-27 (same worker/periods as
# before, just re-located further down the table)
$periodsMartha = @("1801","1712","1711","1710","1709")
for ($i = 0; $i -lt $periodsMartha.Length; $i++) {
  $r = 23 + $i
  $ws.Cells.Item($r, 2).Value = "CC"
  $ws.Cells.Item($r, 3).Value = "1143332154"
  $ws.Cells.Item($r, 4).Value = "MARTHA LIGIA FONNEGRA GEORGE"
  $ws.Cells.Item($r, 5).Value = $periodsMartha[$i]
  $ws.Cells.Item($r, 6).Value = 29600
  $ws.Cells.Item($r, 7).Value = 840000
}

# --- Update the summary header ----------------------------------------------
# Total mora: 7*36400 + 5*29600 = 402800
$ws.Range("E11").Value = 402800
# Cant. Trabajadores: 2
$ws.Range("C13").Value = 2
# Cant. Periodos: 12
$ws.Range("F13").Value = 12

Write-Output "edit complete"
